$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values added in this edit
$ws.Range("J4").Value = 93.65
$ws.Range("J11").Value = 93.65
$ws.Range("E12").Value = 142.8722

# Update the active selection to match the saved view state
$ws.Range("I17").Select() | Out-Null
